# Fix maintenance log form
#
# The "choices" sheet's common_spare_parts list drops nine parts that are no
# longer tracked (burner, fuse, gas_regulator, hose, leak_detector,
# refrigerant, thermocouple, thermometer, voltage_regulator). Removing the
# rows shifts the remaining entries up and makes "choices" the active sheet
# (it was "survey" before), with the selection left on the last remaining
# spare-parts row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

$partsToRemove = @(
    "burner",
    "fuse",
    "gas_regulator",
    "hose",
    "leak_detector",
    "refrigerant",
    "thermocouple",
    "thermometer",
    "voltage_regulator"
)

$lastRow = $ws.UsedRange.Rows.Count

foreach ($partName in $partsToRemove) {
    for ($r = $lastRow; $r -ge 1; $r--) {
        $cellValue = $ws.Cells.Item($r, 2).Value2
        if ($cellValue -eq $partName) {
            $ws.Rows($r).Delete() | Out-Null
            break
        }
    }
}

# "choices" becomes the active/selected sheet, with the cursor on B20
# (previously "survey" was the selected sheet).
$ws.Activate() | Out-Null
$ws.Range("B20").Select() | Out-Null
